# esqlabsR CompiledDataSet.xlsx update
# - Bump a couple of house-keeping identifiers on the workbook (file version /
#   revision GUID) and refresh the saved window geometry.
# - Add a new worksheet "Stevens_2012_placebo" (gastric-emptying data digitised
#   from Stevens 2012) at the end of the workbook and make it the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new worksheet after the existing ones.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Stevens_2012_placebo"

# ---------------------------------------------------------------------------
# 2. Header row - same layout as "TestSheet_1" but with the time/fraction/
#    error columns relabelled for this (fraction-remaining) data set.
# ---------------------------------------------------------------------------
$headers = @(
    "Study Id", "Patient Id", "Organ", "Compartment", "Species", "Gender",
    "Dose [unit]", "Molecule", "MW", "Time [min]", "Fraction [%]", "Error [%]",
    "Route", "Group Id"
)
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# ---------------------------------------------------------------------------
# 3. Data rows. Columns: row | J (Time [min]) | K (Fraction [%]) |
#    L (Error [%], optional) | group-name-id (see $groupNames below)
# ---------------------------------------------------------------------------
$groupNames = @{
    "67" = "Placebo_total"
    "68" = "Sita_total"
    "69" = "Placebo_proximal"
    "70" = "Sita_proximal"
    "71" = "Placebo_distal"
    "72" = "Sita_dist"
}

$data = @"
2|0|100||67
3|13.172268907563|93.75||67
4|29.403361344537799|84.1666666666666||67
5|44.647058823529399|72.5||67
6|73.079831932773104|52.9166666666666||67
7|88.273109243697405|46.25||67
8|105.48319327730999|39.5833333333333||67
9|118.66386554621801|32.5||67
10|149.029411764705|21.25||67
11|180.35294117647001|15||67
12|210.65966386554601|9.5833333333333197||67
13|240.94537815126|6.25||67
14|0|100||68
15|14.218487394957901|90||68
16|27.403361344537799|82.5||68
17|43.634453781512597|72.9166666666666||68
18|57.8403361344537|64.1666666666666||68
19|74.067226890756302|55||68
20|89.260504201680604|48.3333333333333||68
21|105.462184873949|41.6666666666666||68
22|119.642857142857|35.4166666666666||68
23|150.004201680672|24.5833333333333||68
24|179.31932773109199|17.5||68
25|209.62605042016801|12.0833333333333||68
26|240.915966386554|9.1666666666666803||68
27|0.967598097502957|74.482758620689594|4.827586206896612|69
28|16.306480380499401|60.689655172413801|2.0689655172413026|69
29|30.659928656361402|50.689655172413801|3.1034482758620001|69
30|45.026753864447002|41.724137931034399|1.7241379310344982|69
31|61.475921521997499|33.793103448275801|2.7586206896551957|69
32|74.839476813317404|27.241379310344801|2.4137931034482989|69
33|90.303210463733507|23.103448275862|2.4137931034482989|69
34|105.76248513674101|18.620689655172399|2.413793103448203|69
35|121.23513674197299|15.1724137931034|2.758620689655201|69
36|150.133769322235|10|2.7586206896550998|69
37|182.16260404280601|6.8965517241379297|2.068965517241371|69
38|211.09690844233|4.4827586206896601||69
39|242.11355529131899|3.1034482758620601||69
40|0.967598097502957|74.482758620689594|4.827586206896612|70
41|16.306480380499401|60.689655172413801|2.0689655172413026|70
42|30.659928656361402|50.689655172413801|3.1034482758620001|70
43|45.026753864447002|41.724137931034399|1.7241379310344982|70
44|61.475921521997499|33.793103448275801|2.7586206896551957|70
45|74.839476813317404|27.241379310344801|2.4137931034482989|70
46|90.303210463733507|23.103448275862|2.4137931034482989|70
47|105.76248513674101|18.620689655172399|2.413793103448203|70
48|121.23513674197299|15.1724137931034|2.758620689655201|70
49|150.133769322235|10|2.7586206896550998|70
50|182.16260404280601|6.8965517241379297|2.068965517241371|70
51|211.09690844233|4.4827586206896601||70
52|242.11355529131899|3.1034482758620601||70
53|0|25.579399141630901|2.7467811158797986|71
54|14.482758620689699|32.446351931330398|2.9184549356223002|71
55|28.965517241379299|32.103004291845401|2.7467811158798021|71
56|45.517241379310398|30.2145922746781|2.4034334763949019|71
57|58.965517241379402|29.0128755364806|2.2317596566523008|71
58|74.482758620689694|25.7510729613733|2.0600858369097992|71
59|88.965517241379402|22.317596566523601|1.7167381974249025|71
60|104.482758620689|19.399141630901202|1.716738197424803|71
61|120|16.480686695278902|1.5450643776824009|71
62|150|11.158798283261801|1.2017167381974208|71
63|180|8.2403433476394792|0.85836909871243883|71
64|211.03448275861999|5.3218884120171701|1.0300429184549404|71
65|238.96551724137899|4.1201716738197396|0.51502145922746978|71
66|0|24.206008583690899|2.7467811158797986|72
67|14.482758620689699|31.244635193133|2.5751072961373005|72
68|26.8965517241379|32.274678111587903|2.9184549356223037|72
69|44.482758620689701|31.244635193133|2.4034334763949019|72
70|61.034482758620797|29.356223175965599|2.7467811158798021|72
71|74.482758620689694|27.9828326180257|2.2317596566524003|72
72|90.000000000000099|25.0643776824034|2.7467811158797986|72
73|104.482758620689|21.8025751072961|2.7467811158797986|72
74|118.965517241379|19.2274678111588|2.4034334763947989|72
75|148.96551724137899|13.5622317596566|2.2317596566524003|72
76|180|9.6137339055793998|2.4034334763948006|72
77|208.96551724137899|6.5236051502145997|1.8884120171673704|72
78|240|5.1502145922746703|1.8884120171673899|72
"@

$rows = $data -split "`n"
foreach ($rawLine in $rows) {
    $line = $rawLine.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'

    $r = [int]$parts[0]
    $timeVal = [double]$parts[1]
    $fracVal = [double]$parts[2]
    $errText = $parts[3]
    $groupId = $parts[4]

    $ws.Cells.Item($r, 10).Value = $timeVal   # J: Time [min]
    $ws.Cells.Item($r, 11).Value = $fracVal   # K: Fraction [%]
    if ($errText.Length -gt 0) {
        $ws.Cells.Item($r, 12).Value = [double]$errText   # L: Error [%]
    }
    $ws.Cells.Item($r, 14).Value = $groupNames[$groupId]  # N: Group Id
}

# ---------------------------------------------------------------------------
# 4. Make the new sheet the active one (matches the saved workbookView's
#    activeTab pointing at this, the 3rd, sheet).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("I2").Select()
